$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (Price / Volume) to Text format so numeric-looking
# strings (e.g. "34.831.23", "1.01", "0.619") are preserved verbatim as text
# instead of being auto-coerced into numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '34.831.23'
$ws.Range("E2").Value = '  -1.09%  '
$ws.Range("D3").Value = '1.834.00'
$ws.Range("E3").Value = '  +0.79%  '
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").Value = '231.27'
$ws.Range("E5").Value = '  -1.00%  '
$ws.Range("D6").Value = '0.619'
$ws.Range("E6").Value = '  +0.43%  '
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("D8").Value = '39.63'
$ws.Range("E8").Value = '  -5.19%  '
$ws.Range("D9").Value = '0.326'
$ws.Range("E9").Value = '  -0.60%  '
$ws.Range("E10").Value = '  -0.61%  '
$ws.Range("D11").Value = '0.0986'
$ws.Range("E11").Value = '  -1.45%  '
$ws.Range("D12").Value = '2.098.76'
$ws.Range("E12").Value = '  +0.75%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Value = '11.32'
$ws.Range("E13").Value = '  +1.61%  '
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = '0.671'
$ws.Range("E14").Value = '  +1.24%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.810.37'
$ws.Range("E15").Value = '  -0.63%  '
$ws.Range("D16").Value = '4.63'
$ws.Range("E16").Value = '  -1.18%  '
$ws.Range("D17").Value = '34.841.06'
$ws.Range("E17").Value = '  -0.79%  '
$ws.Range("D18").Value = '69.60'
$ws.Range("E18").Value = '  -0.31%  '
$ws.Range("D19").Value = '0.0₃0786'
$ws.Range("E19").Value = '  -1.18%  '
$ws.Range("D20").Value = '240.55'
$ws.Range("D21").Value = '12.16'
$ws.Range("E21").Value = '  +1.67%  '
$ws.Range("D22").Value = '4.68'
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("E23").Value = '  +0.22%  '
$ws.Range("D24").Value = '2.26'
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("D25").Value = '171.38'
$ws.Range("E25").Value = '  -0.92%  '
$ws.Range("D26").Value = '7.75'
$ws.Range("E26").Value = '  -1.68%  '
$ws.Range("E27").Value = '  +2.37%  '
$ws.Range("D28").Value = '17.35'
$ws.Range("E28").Value = '  -1.22%  '
$ws.Range("D29").Value = '1.52'
$ws.Range("E29").Value = '  -6.70%  '
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("D31").Value = '0.0551'
$ws.Range("E31").Value = '  -1.12%  '
$ws.Range("D32").Value = '3.93'
$ws.Range("E32").Value = '  -3.01%  '
$ws.Range("E33").Value = '  -1.83%  '
$ws.Range("D34").Value = '1.85'
$ws.Range("E34").Value = '  +3.44%  '
$ws.Range("E35").Value = '  +7.17%  '
$ws.Range("E36").Value = '  +11.03%  '
$ws.Range("E37").Value = '  +1.36%  '
$ws.Range("D38").Value = '91.29'
$ws.Range("E38").Value = '  -2.47%  '
$ws.Range("D39").Value = '1.06'
$ws.Range("E39").Value = '  +5.73%  '
$ws.Range("D40").Value = '1.339.79'
$ws.Range("E40").Value = '  +2.02%  '
$ws.Range("E41").Value = '  -0.68%  '
$ws.Range("D42").Value = '14.53'
$ws.Range("E42").Value = '  -1.70%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '2.26'
$ws.Range("E43").Value = '  -2.90%  '
$ws.Range("B44").Value = 'HuobiToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D44").Value = '2.41'
$ws.Range("E44").Value = '  -2.04%  '
$ws.Range("D45").Value = '2.76'
$ws.Range("E45").Value = '  -0.46%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '6.23'
$ws.Range("E46").Value = '  -1.84%  '
$ws.Range("B47").Value = 'Kaspa'
$ws.Range("C47").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D47").Value = '0.0521'
$ws.Range("E47").Value = '  +1.61%  '
$ws.Range("D48").Value = '2.013.06'
$ws.Range("E48").Value = '  +0.74%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.0673'
$ws.Range("E49").Value = '  +3.62%  '
$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D50").Value = '1.01'
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("D51").Value = '3.25'
$ws.Range("E51").Value = '  +13.69%  '

# Restore the default cell style so no stray formatting is left behind
# (only the text-forcing number format was needed during the write).
$ws.Range("D2:E51").Style = "Normal"

